$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new product ("BRUFEN 400MG 30 TAB") needs to be inserted, in alphabetical
# order, as the 3rd data row (between "ANTI-COX II 15MG 30 TAB." and
# "EREC 100MG 12 F.C. TABLETS"). That pushes every row below it down by one,
# including the totals row and the footer row.
$ws.Rows.Item(9).Insert()

# Re-create the merged ranges for the newly inserted row 9 (Excel's row
# insert does not automatically recreate the merge on the blank inserted
# row in every environment, so do it explicitly to match the other rows).
$ws.Range("A9:B9").Merge()
$ws.Range("C9:G9").Merge()
$ws.Range("H9:K9").Merge()
$ws.Range("L9:M9").Merge()
$ws.Range("N9:O9").Merge()

# Fill in the data for the new row (serial number, name, balance, min, price,
# sell price, number of transactions). The ratio/price columns are stored as
# TEXT in this report (not numbers), so force a text number format first to
# stop Excel from silently re-typing them as numeric values.
$ws.Cells.Item(9, 1).Value = 3
$ws.Cells.Item(9, 3).Value = "BRUFEN 400MG 30 TAB"

$ws.Cells.Item(9, 8).NumberFormat = "@"
$ws.Cells.Item(9, 8).Value = "1:0"

$ws.Cells.Item(9, 12).NumberFormat = "@"
$ws.Cells.Item(9, 12).Value = "1"

$ws.Cells.Item(9, 14).NumberFormat = "@"
$ws.Cells.Item(9, 14).Value = "78.00"

$ws.Cells.Item(9, 16).NumberFormat = "@"
$ws.Cells.Item(9, 16).Value = "25.7400"

$ws.Cells.Item(9, 17).NumberFormat = "@"
$ws.Cells.Item(9, 17).Value = "0:1"

# Renumber the serial number ("م") column for the rows that followed the
# inserted one (they were shifted down by one row but kept their old
# numbers).
For ($i = 10; $i -le 22; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 6
}

# The grand-total cell (selling price column) moved from row 22 to row 23
# and increases by the new item's selling price (457.24 + 25.74 = 482.98).
$ws.Cells.Item(23, 16).Value = 482.98

# Update the generated timestamp shown in the footer row (now row 24).
$ws.Cells.Item(24, 1).Value = "Friday, 29 August, 2025 4:36 PM"
